$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 351, shifting existing rows 351..391 down to 352..392.
$ws.Rows("351:351").Insert()

# Populate the newly inserted row 351 with the new data.
$ws.Cells.Item(351, 1).Value = 4
$ws.Cells.Item(351, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(351, 3).Value = "Los Lagos"
$ws.Cells.Item(351, 4).Value = 44946
$ws.Cells.Item(351, 5).Value = 10
$ws.Cells.Item(351, 6).Value = 100112040
$ws.Cells.Item(351, 7).Value = "Cilantro"
$ws.Cells.Item(351, 8).Value = "Sin especificar"
$ws.Cells.Item(351, 9).Value = "Primera"
$ws.Cells.Item(351, 10).Value = 180
$ws.Cells.Item(351, 11).Value = 8000
$ws.Cells.Item(351, 12).Value = 8000
$ws.Cells.Item(351, 13).Value = 8000
$ws.Cells.Item(351, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(351, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(351, 16).Value = 4000
$ws.Cells.Item(351, 17).Value = 2
$ws.Cells.Item(351, 18).Value = "Hortaliza"
